# Automatische test-sync: 2025-06-19 22:08:50
# Append a new mail-log entry to the "Logs" sheet and refresh the
# "Dashboard" summary counts accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row being appended to the Logs sheet.
$newRow = 37

$logs.Cells.Item($newRow, 1).Value = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item($newRow, 4).Value = "Klacht / Probleem"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 22:08:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Update the Dashboard summary: the "Klacht / Probleem" count goes up by
# one, and the "IT / Technisch probleem" / "Offerte / Prijsaanvraag" rows
# swap their category labels (counts stay as they were).
$dashboard.Cells.Item(5, 2).Value = 4
$dashboard.Cells.Item(6, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(8, 1).Value = "IT / Technisch probleem"

# The conditional-formatting ranges on the Logs sheet covered the old
# data extent (rows 2-36); stretch them to include the freshly appended
# row 37 as well.
$catRange = $logs.Range("D2:D36").FormatConditions
for ($i = 1; $i -le $catRange.Count; $i++) {
    $catRange.Item($i).ModifyAppliesToRange($logs.Range("D2:D37"))
}

$answeredRange = $logs.Range("G2:G36").FormatConditions
for ($i = 1; $i -le $answeredRange.Count; $i++) {
    $answeredRange.Item($i).ModifyAppliesToRange($logs.Range("G2:G37"))
}
